$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 41
$ws.Range("F4").Value = 17
$ws.Range("F5").Value = 164
$ws.Range("F6").Value = 1039
$ws.Range("F7").Value = 624
$ws.Range("F8").Value = 7928
$ws.Range("F10").Value = 190
$ws.Range("F11").Value = 6821
$ws.Range("F12").Value = 158
$ws.Range("F14").Value = 4851
$ws.Range("F16").Value = 5261
$ws.Range("F18").Value = 309
$ws.Range("F19").Value = 312
$ws.Range("F20").Value = 425
$ws.Range("F21").Value = 307
$ws.Range("F22").Value = 251
$ws.Range("F23").Value = 134
$ws.Range("F26").Value = 8940
$ws.Range("F27").Value = 68
$ws.Range("F28").Value = 1592
$ws.Range("F32").Value = 817
$ws.Range("F33").Value = 70
$ws.Range("F37").Value = 1837
$ws.Range("F38").Value = 233
$ws.Range("F39").Value = 1134
$ws.Range("F41").Value = 4658
$ws.Range("F42").Value = 24
$ws.Range("F44").Value = 62
$ws.Range("F45").Value = 137
$ws.Range("F46").Value = 69
$ws.Range("F47").Value = 901
$ws.Range("F48").Value = 1222
$ws.Range("F49").Value = 55

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 36
$ws.Range("F6").Value = 22
$ws.Range("F13").Value = 4
$ws.Range("F17").Value = 882

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 41
$ws.Range("F5").Value = 17
$ws.Range("F6").Value = 164
$ws.Range("F7").Value = 1039
$ws.Range("F8").Value = 624
$ws.Range("F9").Value = 7928
$ws.Range("F11").Value = 190
$ws.Range("F12").Value = 6821
$ws.Range("F13").Value = 158
$ws.Range("F16").Value = 4851
$ws.Range("F18").Value = 5262
$ws.Range("F20").Value = 309
$ws.Range("F21").Value = 312
$ws.Range("F22").Value = 425
$ws.Range("F23").Value = 307
$ws.Range("F24").Value = 251
$ws.Range("F25").Value = 134
$ws.Range("F27").Value = 8940
$ws.Range("F28").Value = 68
$ws.Range("F29").Value = 1592
$ws.Range("F32").Value = 816
$ws.Range("F33").Value = 70
$ws.Range("F37").Value = 1837
$ws.Range("F38").Value = 233
$ws.Range("F39").Value = 1134
$ws.Range("F41").Value = 4658
$ws.Range("F42").Value = 24
$ws.Range("F44").Value = 62
$ws.Range("F45").Value = 137
$ws.Range("F46").Value = 69
$ws.Range("F47").Value = 901
$ws.Range("F48").Value = 1222
$ws.Range("F49").Value = 55
